$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the gym names in column A (rows 2-11) to the new set of gyms,
# and replace the old trailing caption row (A12) with the newly added gym
# "Dhar Fitness Studio" (new shared string introduced by this commit).
$ws.Range("A2").Value  = "Cross Roads Fitness Studio By Javeed"
$ws.Range("A3").Value  = "Brood Fitness - Semi Personal Training Gym"
$ws.Range("A4").Value  = "6E GYMS"
$ws.Range("A5").Value  = "Gold's Gym"
$ws.Range("A6").Value  = "Phoenix Fitness"
$ws.Range("A7").Value  = "Berlin Fitness"
$ws.Range("A8").Value  = "All About Eve"
$ws.Range("A9").Value  = "Sm Fitness Club"
$ws.Range("A10").Value = "Core Fitness Gym"
$ws.Range("A11").Value = "BLFF -Bruce Lee Freestyle Fitness"
$ws.Range("A12").Value = "Dhar Fitness Studio"
